$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 69570
$ws.Range("E9").Value = 191357970

$ws.Range("C10").Value = 278194
$ws.Range("E10").Value = 1752207808

$ws.Range("C14").Value = 119203
$ws.Range("E14").Value = 379400089

$ws.Range("C104").Value = 22091
$ws.Range("E104").Value = 84771211

$ws.Range("C115").Value = 17550
$ws.Range("E115").Value = 38608142

$ws.Range("C120").Value = 2330
$ws.Range("E120").Value = 4372951

$ws.Range("C164").Value = 50567
$ws.Range("E164").Value = 168401925

$ws.Range("C168").Value = 284943
$ws.Range("D168").Value = 58122
$ws.Range("E168").Value = 1208719126

$ws.Range("C169").Value = 562565
$ws.Range("E169").Value = 1284336204

$ws.Range("C170").Value = 367288
$ws.Range("E170").Value = 2844202045

$ws.Range("C171").Value = 115108
$ws.Range("E171").Value = 444769728

$ws.Range("C173").Value = 54384
$ws.Range("E173").Value = 151848977

$ws.Range("C174").Value = 357168
$ws.Range("E174").Value = 1016719883

$ws.Range("C175").Value = 125514
$ws.Range("E175").Value = 811815500

$ws.Range("C177").Value = 96747
$ws.Range("E177").Value = 174708298

$ws.Range("C179").Value = 235659
$ws.Range("E179").Value = 812286671

$ws.Range("C180").Value = 141456
$ws.Range("E180").Value = 340126057

$ws.Range("C203").Value = 13103
$ws.Range("E203").Value = 33006180

$ws.Range("C204").Value = 4756
$ws.Range("E204").Value = 11726703

$ws.Range("C213").Value = 3634
$ws.Range("E213").Value = 11192315

$ws.Range("C221").Value = 2494
$ws.Range("E221").Value = 7098337

$ws.Range("C240").Value = 5414
$ws.Range("E240").Value = 14226741

$ws.Range("C276").Value = 216631
$ws.Range("E276").Value = 1209907101
